# Generate Report for Handoff
# Updates the localization-status workbook with a fresh handoff report:
#  - Overview sheet: bump "Latest HO Xliff Generate Date" for the
#    096a5302 / a663fc5b / f9e203fd / fbe999d2 rows (rows 4-7).
#  - zh-cn sheet: those same rows move from "low" priority to "ht"
#    (handed-off) and get a refreshed "Latest Handoff Datetime".
#  - de-de sheet: those same rows also move from "low" priority to "ht".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview: refresh the Latest HO Xliff Generate Date for rows 4-7
$overview.Range("G4:G7").Value = "2016-08-18 20:32:16"

# zh-cn: Priority low -> ht, and new Latest Handoff Datetime for rows 4-7
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-18 20:32:09"

# de-de: Priority low -> ht for rows 4-7. The "Latest Handoff Datetime"
# column here shares its text with the Overview "Latest HO Xliff Generate
# Date" column (same underlying value), so refresh it the same way.
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-18 20:32:16"
